$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.487.08"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.507.34"
$ws.Range("E3").Value = "  -3.08%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.81"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.17"
$ws.Range("E6").Value = "  -4.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.504.80"
$ws.Range("E7").Value = "  -3.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.505"
$ws.Range("E9").Value = "  +3.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.70"
$ws.Range("E10").Value = "  -3.17%  "
$ws.Range("E11").Value = "  -5.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.402"
$ws.Range("E12").Value = "  -3.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.096.21"
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("E14").Value = "  -7.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.61"
$ws.Range("E15").Value = "  -4.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.513.36"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.389.40"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.62"
$ws.Range("E19").Value = "  -8.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.08"
$ws.Range("E20").Value = "  -4.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.54"
$ws.Range("E21").Value = "  -3.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "420.75"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("E23").Value = "  -5.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "76.76"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.653.29"
$ws.Range("E25").Value = "  -2.92%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -9.10%  "
$ws.Range("E28").Value = "  -3.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.75"
$ws.Range("E29").Value = "  -8.17%  "
$ws.Range("E30").Value = "  -7.28%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.513.68"
$ws.Range("E32").Value = "  -2.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.153"
$ws.Range("E33").Value = "  -4.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.09"
$ws.Range("E34").Value = "  -5.46%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -10.72%  "
$ws.Range("E37").Value = "  -5.40%  "
$ws.Range("E38").Value = "  -5.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "173.07"
$ws.Range("E39").Value = "  -2.38%  "
$ws.Range("E40").Value = "  -8.99%  "
$ws.Range("E41").Value = "  -7.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.94"
$ws.Range("E42").Value = "  -5.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.849"
$ws.Range("E43").Value = "  -5.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.53"
$ws.Range("E44").Value = "  -1.58%  "
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.34"
$ws.Range("E47").Value = "  -10.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.01"
$ws.Range("E48").Value = "  -2.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.10"
$ws.Range("E49").Value = "  -5.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.77"
$ws.Range("E50").Value = "  -5.03%  "
$ws.Range("E51").Value = "  -7.62%  "
